$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 20.02.2022 01:45"

# Update row 10 (EuroOil Opuštěná): D10 delta and E10 date, both become numeric
$ws.Range("D10").Value = -0.4
$ws.Range("E10").Value = 44612.06721064815
$ws.Range("E10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
